# Weekly fruit/hortaliza update: insert 3 new "Sandia" price rows (one per
# quality grade, dated 2023-10-30 / serial 45229) ahead of the existing
# history, pushing the rest of the table down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 56 - everything that used
# to live at rows 56:85 now lives at 59:88, and the sheet's used range
# grows from A1:R85 to A1:R88 automatically.
$ws.Rows("56:58").Insert()

$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$fecha     = 45229
$codreg    = 15
$catId     = 100112028
$categoria = "Sandia"
$variedad  = "Sin especificar"
$unidad    = "$/kilo (volumen en unidades)"
$origen    = "Perú"
$kgUnid    = 1
$clasif    = "Hortaliza"

$newRows = @(
    @{ Row = 56; Calidad = "Extra";   Volumen = 160; PMin = 450; PMax = 460; PProm = 456 },
    @{ Row = 57; Calidad = "Primera"; Volumen = 650; PMin = 440; PMax = 450; PProm = 445 },
    @{ Row = 58; Calidad = "Segunda"; Volumen = 350; PMin = 440; PMax = 450; PProm = 446 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = 1
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $catId
    $ws.Cells.Item($row, 7).Value  = $categoria
    $ws.Cells.Item($row, 8).Value  = $variedad
    $ws.Cells.Item($row, 9).Value  = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $origen
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $kgUnid
    $ws.Cells.Item($row, 18).Value = $clasif
}
